$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# New task/appointment rows 104-108 (GUI for renting and appointments).
# Row 104 copies the "header-ish" look of row 103 (A46/B68/C69/D66 styles);
# rows 105-108 copy the plain data-row look used by rows 96/98/100
# (A64/B64/C70/D64 styles).
# ---------------------------------------------------------------------------

$ws.Range("A103:D103").Copy()
$ws.Range("A104:D104").PasteSpecial(-4122)

$ws.Range("A100:D100").Copy()
$ws.Range("A105:D108").PasteSpecial(-4122)

$ws.Application.CutCopyMode = $false

$ws.Cells.Item(104,1).Value = "RFID"
$ws.Cells.Item(104,2).Value = "Ilia"
$ws.Cells.Item(104,3).Value = "12.05. || 16:00"
$ws.Cells.Item(104,4).Value = 180

$ws.Cells.Item(105,1).Value = "Website functionality"
$ws.Cells.Item(105,2).Value = "Angel"
$ws.Cells.Item(105,3).Value = "12.05. || 19:00"
$ws.Cells.Item(105,4).Value = 60

$ws.Cells.Item(106,1).Value = "Website functionality"
$ws.Cells.Item(106,2).Value = "Angel Georgi and Ilia"
$ws.Cells.Item(106,3).Value = "12.05. || 20:00"
$ws.Cells.Item(106,4).Value = 120

$ws.Cells.Item(107,1).Value = "Process report"
$ws.Cells.Item(107,2).Value = "Georgi"
$ws.Cells.Item(107,3).Value = "14.05. || 16:00"
$ws.Cells.Item(107,4).Value = 60

$ws.Cells.Item(108,1).Value = "Database to class conversion"
$ws.Cells.Item(108,2).Value = "Ilia"
$ws.Cells.Item(108,3).Value = "14.05. || 21:00"
$ws.Cells.Item(108,4).Value = 180

# ---------------------------------------------------------------------------
# Append a new blank row 112 (same look as the existing trailing blank rows),
# pushing the used range down to A1:S112.
# ---------------------------------------------------------------------------

$ws.Range("A111:D111").Copy()
$ws.Range("A112:D112").PasteSpecial(-4122)
$ws.Rows.Item(112).RowHeight = $ws.Rows.Item(111).RowHeight
$ws.Application.CutCopyMode = $false

# ---------------------------------------------------------------------------
# Update the view so it matches where the author ended up after editing.
# ---------------------------------------------------------------------------

$ws.Range("B106").Select()
